{"js": "// \"turn summary to abstract\"\n//\n// The document currently has:\n//   Heading2 \"Summary\"\n//   (empty paragraph)\n//   paragraph with the long italic justification text (direct run formatting, no paragraph style)\n//   (empty paragraph, still carrying the italic paragraph-mark formatting)\n//   Heading2 \"Background\"\n//   ...\n//\n// We need to:\n//   1. Drop the \"Summary\" Heading2 paragraph and the empty paragraph right after it.\n//   2. Turn the long-text paragraph into the built-in \"Abstract\" paragraph style\n//      (which also clears the manual italic formatting that the runs used to carry,\n//      since the text is no longer meant to look like an italic pull-quote but a\n//      normal abstract paragraph).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Locate the \"Summary\" Heading 2 paragraph robustly (don't assume a fixed index).\nlet summaryIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.style === \"Heading 2\" && para.text.trim() === \"Summary\") {\n    summaryIndex = i;\n    break;\n  }\n}\n\nif (summaryIndex === -1) {\n  throw new Error('Could not find the \"Summary\" Heading 2 paragraph.');\n}\n\nconst summaryParagraph = paragraphs.items[summaryIndex];\nconst blankParagraph = paragraphs.items[summaryIndex + 1];\nconst abstractParagraph = paragraphs.items[summaryIndex + 2];\n\n// Re-style the long text paragraph as the Abstract style first (while references\n// are still valid), then remove the old heading + the blank line that followed it.\nabstractParagraph.style = \"Abstract\";\nawait context.sync();\n\nsummaryParagraph.delete();\nblankParagraph.delete();\nawait context.sync();\n", "ps1": "# \"turn summary to abstract\"\n#\n# The document currently has:\n#   Heading2 \"Summary\"\n#   (empty paragraph)\n#   paragraph with the long italic justification text (direct run formatting, no paragraph style)\n#   (empty paragraph, still carrying the italic paragraph-mark formatting)\n#   Heading2 \"Background\"\n#   ...\n#\n# We need to:\n#   1. Drop the \"Summary\" Heading2 paragraph and the empty paragraph right after it.\n#   2. Turn the long-text paragraph into the built-in \"Abstract\" paragraph style\n#      (which also clears the manual italic formatting that the runs used to carry).\n\n$d = $word.ActiveDocument\n\n# Locate the \"Summary\" Heading 2 paragraph robustly (don't assume a fixed index).\n$summaryIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  $txt = $p.Range.Text.Trim()\n  if ($p.Style.NameLocal -eq \"Heading 2\" -and $txt -eq \"Summary\") {\n    $summaryIndex = $i\n    break\n  }\n}\n\nif ($summaryIndex -eq -1) {\n  throw \"Could not find the 'Summary' Heading 2 paragraph.\"\n}\n\n# Re-style the long text paragraph (summaryIndex + 2) as Abstract while the\n# index is still valid.\n$abstractParagraph = $d.Paragraphs.Item($summaryIndex + 2)\n$abstractParagraph.Style = \"Abstract\"\n\n# Delete bottom-up with fresh lookups so earlier deletions don't invalidate\n# later paragraph references: first the blank paragraph, then the heading.\n$d.Paragraphs.Item($summaryIndex + 1).Range.Delete()\n$d.Paragraphs.Item($summaryIndex).Range.Delete()\n"}
